$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq 45204) {
        $cell.Value = 45207
    }
}
